$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price entry was inserted above the old row 22, shifting the
# following rows (old 22-24) down to 23-25.
$ws.Rows.Item(22).Insert()

$ws.Range("A22").Value = 11
$ws.Range("B22").Value = "Vega Monumental Concepción"
$ws.Range("C22").Value = "Bíobío"
$ws.Range("D22").Value = 44461
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = 100112031
$ws.Range("G22").Value = "Poroto verde"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 33000
$ws.Range("L22").Value = 34000
$ws.Range("M22").Value = 33500
$ws.Range("N22").Value = '$/malla 25 kilos'
$ws.Range("O22").Value = "Región de Arica y Parinacota"
$ws.Range("P22").Value = 1340
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
